$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.654.59"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -1.36%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.594.15"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -1.58%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.01%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'212.05"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -1.23%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'0.513"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -0.91%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  -0.04%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  -1.54%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  -2.86%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'19.64"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -2.01%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  -1.71%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'1.817.28"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -1.67%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'1.605.01"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -1.01%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "'  -2.73%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = "'  -3.09%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'65.33"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +0.88%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'26.626.23"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = "'0.0₃0730"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -2.82%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'208.99"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -2.60%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  +0.09%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'  -2.38%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  -2.31%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  -3.25%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  -1.83%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'145.86"
$ws.Range("D25").Style = "Normal"
$ws.Range("D27").Value = "'7.17"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -3.53%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'0.115"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -1.69%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'15.33"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -1.55%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'0.0505"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -2.35%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  -1.68%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  -3.99%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'0.669"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -11.64%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  -3.46%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'1.301.42"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -3.24%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E37").Value = "'  -4.95%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  -4.25%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  -2.67%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  +0.09%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.792"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -1.19%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'5.35"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +0.51%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'2.19"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -1.75%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'63.27"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -3.12%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'1.730.15"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -1.55%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'89.16"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -0.91%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'1.61"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -1.82%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.803"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -8.52%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.0984"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -2.15%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  -2.61%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'7.56"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -1.35%  "
$ws.Range("E51").Style = "Normal"
